$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "56.720.80"
Set-TextValue $ws.Range("E2") "  -2.70%  "

Set-TextValue $ws.Range("D3") "2.988.52"
Set-TextValue $ws.Range("E3") "  -4.77%  "

Set-TextValue $ws.Range("E4") "  -0.04%  "

Set-TextValue $ws.Range("D5") "498.17"
Set-TextValue $ws.Range("E5") "  -5.33%  "

Set-TextValue $ws.Range("D6") "135.10"
Set-TextValue $ws.Range("E6") "  +1.09%  "

Set-TextValue $ws.Range("E7") "  -0.04%  "

Set-TextValue $ws.Range("D8") "2.984.04"
Set-TextValue $ws.Range("E8") "  -4.91%  "

Set-TextValue $ws.Range("E9") "  -4.31%  "

Set-TextValue $ws.Range("D10") "7.27"
Set-TextValue $ws.Range("E10") "  -0.50%  "

Set-TextValue $ws.Range("D11") "0.106"
Set-TextValue $ws.Range("E11") "  -2.37%  "

Set-TextValue $ws.Range("E12") "  -6.25%  "

Set-TextValue $ws.Range("E13") "  -1.16%  "

Set-TextValue $ws.Range("D14") "3.495.09"
Set-TextValue $ws.Range("E14") "  -5.12%  "

Set-TextValue $ws.Range("D15") "25.40"
Set-TextValue $ws.Range("E15") "  +0.79%  "

Set-TextValue $ws.Range("D16") "56.660.23"
Set-TextValue $ws.Range("E16") "  -2.85%  "

Set-TextValue $ws.Range("D17") "0.0000149"
Set-TextValue $ws.Range("E17") "  -1.92%  "

Set-TextValue $ws.Range("D18") "2.983.82"
Set-TextValue $ws.Range("E18") "  -5.21%  "

Set-TextValue $ws.Range("D19") "5.73"
Set-TextValue $ws.Range("E19") "  -0.33%  "

Set-TextValue $ws.Range("D20") "12.36"
Set-TextValue $ws.Range("E20") "  -4.78%  "

Set-TextValue $ws.Range("E21") "  -1.15%  "

Set-TextValue $ws.Range("D22") "326.83"
Set-TextValue $ws.Range("E22") "  -4.48%  "

Set-TextValue $ws.Range("D23") "1.00"
Set-TextValue $ws.Range("E23") "  -0.03%  "

Set-TextValue $ws.Range("D24") "0.471"
Set-TextValue $ws.Range("E24") "  -7.45%  "

Set-TextValue $ws.Range("D25") "61.84"
Set-TextValue $ws.Range("E25") "  -8.28%  "

Set-TextValue $ws.Range("D26") "0.998"
Set-TextValue $ws.Range("E26") "  -0.17%  "

Set-TextValue $ws.Range("E27") "  -5.11%  "

Set-TextValue $ws.Range("D28") "0.0₃0901"
Set-TextValue $ws.Range("E28") "  -3.45%  "

Set-TextValue $ws.Range("E29") "  +0.08%  "

Set-TextValue $ws.Range("E30") "  -4.54%  "

Set-TextValue $ws.Range("D31") "6.78"
Set-TextValue $ws.Range("E31") "  -1.24%  "

Set-TextValue $ws.Range("D32") "20.61"
Set-TextValue $ws.Range("E32") "  -3.97%  "

Set-TextValue $ws.Range("D33") "1.18"
Set-TextValue $ws.Range("E33") "  -8.12%  "

Set-TextValue $ws.Range("D34") "1.74"
Set-TextValue $ws.Range("E34") "  -6.93%  "

Set-TextValue $ws.Range("D35") "152.43"
Set-TextValue $ws.Range("E35") "  -4.48%  "

Set-TextValue $ws.Range("D36") "4.46"
Set-TextValue $ws.Range("E36") "  -7.56%  "

Set-TextValue $ws.Range("E37") "  -7.16%  "

Set-TextValue $ws.Range("D38") "5.61"
Set-TextValue $ws.Range("E38") "  -9.69%  "

Set-TextValue $ws.Range("D39") "0.0673"
Set-TextValue $ws.Range("E39") "  -1.73%  "

Set-TextValue $ws.Range("D40") "23.31"
Set-TextValue $ws.Range("E40") "  -0.56%  "

Set-TextValue $ws.Range("D41") "3.017.33"
Set-TextValue $ws.Range("E41") "  -4.87%  "

Set-TextValue $ws.Range("E42") "  -9.71%  "

Set-TextValue $ws.Range("E43") "  -0.05%  "

Set-TextValue $ws.Range("B44") "ONDO"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D44") "1.00"
Set-TextValue $ws.Range("E44") "  -7.22%  "

Set-TextValue $ws.Range("B45") "Mantle"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D45") "0.638"
Set-TextValue $ws.Range("E45") "  -7.71%  "

Set-TextValue $ws.Range("D46") "2.221.32"
Set-TextValue $ws.Range("E46") "  -2.81%  "

Set-TextValue $ws.Range("E47") "  -3.44%  "

Set-TextValue $ws.Range("E48") "  -9.25%  "

Set-TextValue $ws.Range("D49") "1.94"
Set-TextValue $ws.Range("E49") "  +6.42%  "

Set-TextValue $ws.Range("E50") "  +0.94%  "

Set-TextValue $ws.Range("D51") "5.75"
Set-TextValue $ws.Range("E51") "  -6.69%  "
